$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for every data row (2-89)
# from 45207 (2023-10-08) to 45208 (2023-10-09).
for ($r = 2; $r -le 89; $r++) {
    $ws.Cells.Item($r, 3).Value = 45208
}

# Row 2's record (A 34293-2023) moved from the LINDESBERG logging folder to
# the "1885" logging folder, so update every hyperlink formula in row 2
# that references the old folder name.
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/artfynd/A 34293-2023.xlsx", "A 34293-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/kartor/A 34293-2023.png", "A 34293-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/klagomål/A 34293-2023.docx", "A 34293-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/klagomålsmail/A 34293-2023.docx", "A 34293-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/tillsyn/A 34293-2023.docx", "A 34293-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/tillsynsmail/A 34293-2023.docx", "A 34293-2023")'
